$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# 1. Relabel the last data point (row 19, "d1") as "g2". Doing this first while the
#    shared string "d1" is only referenced once lets the engine reclaim that slot in
#    place instead of appending a brand-new string at the end of the table.
$ws.Range("D19").Value = "g2"

# 2. A new (theta=46) measurement is inserted near the top of that group: insert a row
#    at 13 and fill in its values.
$ws.Rows.Item(13).Insert()
$ws.Range("A13").Value = 46
$ws.Range("B13").Value = 1320
$ws.Range("D13").Value = "g1"

# 3. Another new (theta=46) measurement is inserted further down (originally between
#    2017 and 2970): insert a row at 16 and fill in its values.
$ws.Rows.Item(16).Insert()
$ws.Range("A16").Value = 46
$ws.Range("B16").Value = 1930
$ws.Range("D16").Value = "a3"

# 4. Relabel row 9's peak from "d" to "a3/d" (a new shared string).
$ws.Range("D9").Value = "a3/d"

# 5. Refresh column C's conversion formula across the whole (now longer) table so it
#    is stored as a single shared formula, same as Excel would do on a fill-down.
$ws.Range("C2:C21").Formula = "=COS(RADIANS(A2))*B2"

# 6. Right-align the peak-label column for all data rows.
$ws.Range("D2:D21").HorizontalAlignment = -4152

# 7. Widen columns B:D to fit the (now longer) numbers/labels.
$ws.Columns.Item(2).ColumnWidth = 11.592447916666666
$ws.Columns.Item(3).ColumnWidth = 12.307291666666666
$ws.Columns.Item(4).ColumnWidth = 4.736979166666667

# 8. Leave the selection on D10, matching where editing finished.
$ws.Range("D10").Select() | Out-Null
